$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = "2026-02-19 16:24:33"
$ws.Cells.Item(2, 11).Value = "7.1 MJ/m2"
$ws.Cells.Item(2, 15).Value = "-0.8 °C"

$ws.Cells.Item(3, 5).Value = "2026-02-19 16:24:36"
$ws.Cells.Item(3, 8).Value = "94%"
$ws.Cells.Item(3, 11).Value = "7.7 MJ/m2"
$ws.Cells.Item(3, 15).Value = "-6.3 °C"

$ws.Cells.Item(4, 5).Value = "2026-02-19 16:24:38"
$ws.Cells.Item(4, 8).Value = "60%"
$ws.Cells.Item(4, 11).Value = "10.7 MJ/m2"
$ws.Cells.Item(4, 15).Value = "11.3 °C"

$ws.Cells.Item(5, 5).Value = "2026-02-19 16:24:41"
$ws.Cells.Item(5, 11).Value = "6.5 MJ/m2"
$ws.Cells.Item(5, 15).Value = "-6.2 °C"

$ws.Cells.Item(6, 5).Value = "2026-02-19 16:24:44"
$ws.Cells.Item(6, 8).Value = "75%"
$ws.Cells.Item(6, 11).Value = "11.3 MJ/m2"
$ws.Cells.Item(6, 13).Value = "15.8 °C 15:40 TU"
$ws.Cells.Item(6, 15).Value = "10.0 °C"

$ws.Cells.Item(7, 5).Value = "2026-02-19 16:24:47"
$ws.Cells.Item(7, 10).Value = "1010.0 hPa"
$ws.Cells.Item(7, 11).Value = "13.2 MJ/m2"

$ws.Cells.Item(8, 5).Value = "2026-02-19 16:24:49"
$ws.Cells.Item(8, 8).Value = "55%"
$ws.Cells.Item(8, 10).Value = "1009.7 hPa"
$ws.Cells.Item(8, 11).Value = "12.8 MJ/m2"
$ws.Cells.Item(8, 15).Value = "9.9 °C"

$ws.Cells.Item(9, 5).Value = "2026-02-19 16:24:52"
$ws.Cells.Item(9, 8).Value = "79%"
$ws.Cells.Item(9, 11).Value = "10.0 MJ/m2"
$ws.Cells.Item(9, 15).Value = "10.3 °C"

$ws.Cells.Item(10, 5).Value = "2026-02-19 16:24:55"
$ws.Cells.Item(10, 8).Value = "68%"
$ws.Cells.Item(10, 11).Value = "11.7 MJ/m2"
$ws.Cells.Item(10, 15).Value = "10.9 °C"

$ws.Cells.Item(11, 5).Value = "2026-02-19 16:24:58"
$ws.Cells.Item(11, 8).Value = "72%"
$ws.Cells.Item(11, 15).Value = "4.8 °C"

$ws.Cells.Item(12, 5).Value = "2026-02-19 16:25:00"
$ws.Cells.Item(12, 15).Value = "10.8 °C"

$ws.Cells.Item(13, 5).Value = "2026-02-19 16:25:03"
$ws.Cells.Item(13, 8).Value = "68%"
$ws.Cells.Item(13, 11).Value = "13.2 MJ/m2"
$ws.Cells.Item(13, 15).Value = "3.8 °C"

$ws.Cells.Item(14, 5).Value = "2026-02-19 16:25:05"
$ws.Cells.Item(14, 8).Value = "48%"
$ws.Cells.Item(14, 11).Value = "12.1 MJ/m2"
$ws.Cells.Item(14, 15).Value = "13.2 °C"

$ws.Cells.Item(15, 5).Value = "2026-02-19 16:25:08"
$ws.Cells.Item(15, 8).Value = "75%"
$ws.Cells.Item(15, 15).Value = "9.9 °C"

$ws.Cells.Item(16, 5).Value = "2026-02-19 16:25:10"
$ws.Cells.Item(16, 8).Value = "74%"
$ws.Cells.Item(16, 11).Value = "11.9 MJ/m2"

$ws.Cells.Item(17, 5).Value = "2026-02-19 16:25:13"
$ws.Cells.Item(17, 11).Value = "11.3 MJ/m2"
$ws.Cells.Item(17, 15).Value = "0.2 °C"

$ws.Cells.Item(18, 5).Value = "2026-02-19 16:25:16"
$ws.Cells.Item(18, 8).Value = "62%"
$ws.Cells.Item(18, 11).Value = "10.5 MJ/m2"
$ws.Cells.Item(18, 15).Value = "11.4 °C"

$ws.Cells.Item(19, 5).Value = "2026-02-19 16:25:18"
$ws.Cells.Item(19, 11).Value = "9.8 MJ/m2"
$ws.Cells.Item(19, 15).Value = "5.2 °C"

$ws.Cells.Item(20, 5).Value = "2026-02-19 16:25:21"
$ws.Cells.Item(20, 8).Value = "91%"
$ws.Cells.Item(20, 11).Value = "12.3 MJ/m2"

$ws.Cells.Item(21, 5).Value = "2026-02-19 16:25:24"
$ws.Cells.Item(21, 8).Value = "68%"
$ws.Cells.Item(21, 10).Value = "1010.1 hPa"
$ws.Cells.Item(21, 11).Value = "12.5 MJ/m2"
$ws.Cells.Item(21, 15).Value = "6.0 °C"

$ws.Cells.Item(22, 5).Value = "2026-02-19 16:25:26"
$ws.Cells.Item(22, 11).Value = "14.9 MJ/m2"

$ws.Cells.Item(23, 5).Value = "2026-02-19 16:25:29"
$ws.Cells.Item(23, 9).Value = "5.9 mm"
$ws.Cells.Item(23, 11).Value = "10.9 MJ/m2"
$ws.Cells.Item(23, 15).Value = "-6.6 °C"

$ws.Cells.Item(24, 5).Value = "2026-02-19 16:25:32"
$ws.Cells.Item(24, 8).Value = "65%"
$ws.Cells.Item(24, 10).Value = "1013.7 hPa"
$ws.Cells.Item(24, 11).Value = "14.4 MJ/m2"
$ws.Cells.Item(24, 15).Value = "8.8 °C"

$ws.Cells.Item(25, 5).Value = "2026-02-19 16:25:35"
$ws.Cells.Item(25, 11).Value = "14.1 MJ/m2"
$ws.Cells.Item(25, 15).Value = "-4.4 °C"

$ws.Cells.Item(26, 5).Value = "2026-02-19 16:25:38"
$ws.Cells.Item(26, 10).Value = "1009.1 hPa"
$ws.Cells.Item(26, 11).Value = "8.4 MJ/m2"
$ws.Cells.Item(26, 15).Value = "2.7 °C"

$ws.Cells.Item(27, 5).Value = "2026-02-19 16:25:40"
$ws.Cells.Item(27, 11).Value = "12.8 MJ/m2"
$ws.Cells.Item(27, 15).Value = "-3.9 °C"

$ws.Cells.Item(28, 5).Value = "2026-02-19 16:25:43"
$ws.Cells.Item(28, 10).Value = "1009.1 hPa"
$ws.Cells.Item(28, 11).Value = "10.9 MJ/m2"
$ws.Cells.Item(28, 15).Value = "8.8 °C"

$ws.Cells.Item(29, 5).Value = "2026-02-19 16:25:46"
$ws.Cells.Item(29, 8).Value = "74%"
$ws.Cells.Item(29, 11).Value = "12.0 MJ/m2"
$ws.Cells.Item(29, 15).Value = "10.9 °C"

$ws.Cells.Item(30, 5).Value = "2026-02-19 16:25:48"
$ws.Cells.Item(30, 8).Value = "79%"
$ws.Cells.Item(30, 10).Value = "1009.2 hPa"
$ws.Cells.Item(30, 11).Value = "9.4 MJ/m2"
$ws.Cells.Item(30, 15).Value = "9.9 °C"

$ws.Cells.Item(31, 5).Value = "2026-02-19 16:25:51"
$ws.Cells.Item(31, 8).Value = "54%"
$ws.Cells.Item(31, 10).Value = "1008.6 hPa"
$ws.Cells.Item(31, 11).Value = "8.5 MJ/m2"

$ws.Cells.Item(32, 5).Value = "2026-02-19 16:25:53"
$ws.Cells.Item(32, 11).Value = "13.7 MJ/m2"

$ws.Cells.Item(33, 5).Value = "2026-02-19 16:25:56"
$ws.Cells.Item(33, 8).Value = "63%"
$ws.Cells.Item(33, 10).Value = "1010.0 hPa"
$ws.Cells.Item(33, 11).Value = "6.7 MJ/m2"
$ws.Cells.Item(33, 15).Value = "3.2 °C"

$ws.Cells.Item(34, 5).Value = "2026-02-19 16:25:59"
$ws.Cells.Item(34, 11).Value = "12.8 MJ/m2"
$ws.Cells.Item(34, 12).Value = "70.9 km/h - 268º 15:56 TU"
$ws.Cells.Item(34, 15).Value = "-2.1 °C"

$ws.Cells.Item(35, 5).Value = "2026-02-19 16:26:02"
$ws.Cells.Item(35, 10).Value = "1015.1 hPa"
$ws.Cells.Item(35, 11).Value = "14.5 MJ/m2"

$ws.Cells.Item(36, 5).Value = "2026-02-19 16:26:04"
$ws.Cells.Item(36, 10).Value = "1009.5 hPa"
$ws.Cells.Item(36, 11).Value = "10.4 MJ/m2"
$ws.Cells.Item(36, 15).Value = "11.8 °C"

$ws.Cells.Item(37, 5).Value = "2026-02-19 16:26:07"
$ws.Cells.Item(37, 8).Value = "77%"
$ws.Cells.Item(37, 13).Value = "11.9 °C 15:49 TU"
$ws.Cells.Item(37, 15).Value = "5.2 °C"

$ws.Cells.Item(38, 5).Value = "2026-02-19 16:26:10"
$ws.Cells.Item(38, 11).Value = "11.8 MJ/m2"
$ws.Cells.Item(38, 15).Value = "11.5 °C"

$ws.Cells.Item(39, 5).Value = "2026-02-19 16:26:12"
$ws.Cells.Item(39, 11).Value = "13.8 MJ/m2"
$ws.Cells.Item(39, 15).Value = "-6.1 °C"

$ws.Cells.Item(40, 5).Value = "2026-02-19 16:26:15"
$ws.Cells.Item(40, 8).Value = "80%"
$ws.Cells.Item(40, 15).Value = "5.6 °C"

$ws.Cells.Item(41, 5).Value = "2026-02-19 16:26:17"
$ws.Cells.Item(41, 11).Value = "15.0 MJ/m2"
$ws.Cells.Item(41, 15).Value = "14.1 °C"

$ws.Cells.Item(42, 5).Value = "2026-02-19 16:26:20"
$ws.Cells.Item(42, 15).Value = "11.3 °C"

$ws.Cells.Item(43, 5).Value = "2026-02-19 16:26:23"
$ws.Cells.Item(43, 11).Value = "13.1 MJ/m2"
$ws.Cells.Item(43, 15).Value = "9.0 °C"

$ws.Cells.Item(44, 5).Value = "2026-02-19 16:26:25"
$ws.Cells.Item(44, 11).Value = "9.8 MJ/m2"

$ws.Cells.Item(45, 5).Value = "2026-02-19 16:26:28"
$ws.Cells.Item(45, 8).Value = "88%"
$ws.Cells.Item(45, 10).Value = "1014.2 hPa"
$ws.Cells.Item(45, 11).Value = "7.5 MJ/m2"
$ws.Cells.Item(45, 15).Value = "2.2 °C"

$ws.Cells.Item(46, 5).Value = "2026-02-19 16:26:31"
$ws.Cells.Item(46, 10).Value = "1014.5 hPa"
$ws.Cells.Item(46, 11).Value = "14.6 MJ/m2"
$ws.Cells.Item(46, 15).Value = "12.7 °C"
